$d = $word.ActiveDocument

# Locate the "TABLA DE CONTENIDO" heading paragraph.
$tocIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -match "TABLA DE CONTENIDO") {
        $tocIndex = $i
        break
    }
}

if ($tocIndex -ne -1) {
    # The bullet list right under the heading uses the "List Paragraph" style;
    # walk forward while that style continues to include the whole list.
    $lastIndex = $tocIndex
    $j = $tocIndex + 1
    while ($j -le $d.Paragraphs.Count -and $d.Paragraphs.Item($j).Style.NameLocal -eq "List Paragraph") {
        $lastIndex = $j
        $j = $j + 1
    }

    # "cambio de tamaño letra": bump the heading + its list items to 28pt
    # (sz/szCs = 56 half-points in OOXML).
    for ($i = $tocIndex; $i -le $lastIndex; $i++) {
        $p = $d.Paragraphs.Item($i)
        $p.Range.Font.Size = 28
        $p.Range.Font.SizeBi = 28
    }
}

# Repagination side-effect: the enlarged "TABLA DE CONTENIDO" text pushes the
# page break earlier, so the stale <w:lastRenderedPageBreak/> cached on the
# "INTRODUCCION" run no longer applies there. Touch that run (without
# actually changing its text) so the layout cache is refreshed and the stale
# marker is dropped.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "INTRODUCCION") {
        $r = $p.Range
        $r2 = $d.Range($r.Start, $r.End - 1)
        $r2.InsertBefore("")
        break
    }
}
